# "filtrado por mes y punto" - replace the single R1 sample row with the
# R3 point measurements for the month, adding rows 3-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: Fecha(serial), Punto, Hora Inicio, Hora Fin, Duracion, Tiempo estab., LA,F,eq, LA,F,10, LA,F,90, Estandar
$data = @(
    @(45014, "R3", "14:35", "14:50", 15, 1, 64.7,               65.59999999999999, 63.1, 75),
    @(45012, "R3", "20:17", "20:32", 15, 1, 66.59999999999999,  69.5,              62.8, 75),
    @(45007, "R3", "12:58", "13:13", 15, 2, 63.6,               64.7,              61.7, 75),
    @(45006, "R3", "19:47", "20:02", 15, 1, 65.09999999999999,  66.09999999999999, 63.2, 75),
    @(45002, "R3", "08:12", "08:27", 15, 1, 67.2,               68.8,              65.3, 75),
    @(45001, "R3", "11:22", "11:37", 15, 2, 69.7,               71.7,              66.8, 75)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $ws.Cells.Item($row, 9).Value = $r[8]
    $ws.Cells.Item($row, 10).Value = $r[9]
    $row = $row + 1
}
